$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-10-22 Sunday" "2023-10-23 Monday"

Replace-Text "68×82=" "69×46="
Replace-Text "50×74=" "24×63="
Replace-Text "64×42=" "99×20="
Replace-Text "34×51=" "24×95="
Replace-Text "51×91=" "70×31="

Replace-Text "76×27=" "39×40="
Replace-Text "83×72=" "56×47="
Replace-Text "41×14=" "97×81="
Replace-Text "17×35=" "14×39="
Replace-Text "16×46=" "51×23="

Replace-Text "21×41=" "55×69="
Replace-Text "25×76=" "32×12="
Replace-Text "87×69=" "81×77="
Replace-Text "71×60=" "26×95="
Replace-Text "81×75=" "74×56="

Replace-Text "72×85=" "82×50="
Replace-Text "36×50=" "19×40="
Replace-Text "26×32=" "38×72="
Replace-Text "49×67=" "52×89="
Replace-Text "39×67=" "27×72="

Replace-Text "53×69=" "65×49="
Replace-Text "13×18=" "96×56="
Replace-Text "55×52=" "66×44="
Replace-Text "55×96=" "64×40="
Replace-Text "56×89=" "23×89="
